$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 88.72291666666666
$ws.Range("H2").Value = 266.16875
$ws.Range("I2").Value = 0.7675060578750151
$ws.Range("J2").Value = 0.7675060578750152
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 16.535604
$ws.Range("N2").Value = 49.606812
$ws.Range("O2").Value = 0.2120453146491552
$ws.Range("P2").Value = 0.2120453146491552
$ws.Range("Q2").Value = 1467.087015725
$ws.Range("R2").Value = 13203.783141525
$ws.Range("S2").Value = 0.1627460635372403
$ws.Range("T2").Value = 0.1627460635372404

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 88.72291666666666
$ws.Range("H3").Value = 266.16875
$ws.Range("I3").Value = 0.7675060578750151
$ws.Range("J3").Value = 0.7675060578750152
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 40.62063066666667
$ws.Range("N3").Value = 121.861892
$ws.Range("O3").Value = 0.5209011059384622
$ws.Range("P3").Value = 0.5209011059384622
$ws.Range("Q3").Value = 3603.980829586111
$ws.Range("R3").Value = 32435.827466275
$ws.Range("S3").Value = 0.3997947543615648
$ws.Range("T3").Value = 0.3997947543615648

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 88.72291666666666
$ws.Range("H4").Value = 266.16875
$ws.Range("I4").Value = 0.7675060578750151
$ws.Range("J4").Value = 0.7675060578750152
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 20.825229
$ws.Range("N4").Value = 62.475687
$ws.Range("O4").Value = 0.2670535794123827
$ws.Range("P4").Value = 0.2670535794123827
$ws.Range("Q4").Value = 1847.67505713125
$ws.Range("R4").Value = 16629.07551418125
$ws.Range("S4").Value = 0.2049652399762101
$ws.Range("T4").Value = 0.2049652399762102

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 17.91585
$ws.Range("H5").Value = 53.74755
$ws.Range("I5").Value = 0.1549827702197958
$ws.Range("J5").Value = 0.1549827702197958
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 16.535604
$ws.Range("N5").Value = 49.606812
$ws.Range("O5").Value = 0.2120453146491552
$ws.Range("P5").Value = 0.2120453146491552
$ws.Range("Q5").Value = 296.2494009234
$ws.Range("R5").Value = 2666.2446083106
$ws.Range("S5").Value = 0.03286337027645433
$ws.Range("T5").Value = 0.03286337027645433

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.91585
$ws.Range("H6").Value = 53.74755
$ws.Range("I6").Value = 0.1549827702197958
$ws.Range("J6").Value = 0.1549827702197958
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 40.62063066666667
$ws.Range("N6").Value = 121.861892
$ws.Range("O6").Value = 0.5209011059384622
$ws.Range("P6").Value = 0.5209011059384622
$ws.Range("Q6").Value = 727.7531259294002
$ws.Range("R6").Value = 6549.778133364601
$ws.Range("S6").Value = 0.0807306964088982
$ws.Range("T6").Value = 0.08073069640889821

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 17.91585
$ws.Range("H7").Value = 53.74755
$ws.Range("I7").Value = 0.1549827702197958
$ws.Range("J7").Value = 0.1549827702197958
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 20.825229
$ws.Range("N7").Value = 62.475687
$ws.Range("O7").Value = 0.2670535794123827
$ws.Range("P7").Value = 0.2670535794123827
$ws.Range("Q7").Value = 373.10167897965
$ws.Range("R7").Value = 3357.91511081685
$ws.Range("S7").Value = 0.04138870353444329
$ws.Range("T7").Value = 0.0413887035344433

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.960212333333333
$ws.Range("H8").Value = 26.880637
$ws.Range("I8").Value = 0.077511171905189
$ws.Range("J8").Value = 0.07751117190518901
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 16.535604
$ws.Range("N8").Value = 49.606812
$ws.Range("O8").Value = 0.2120453146491552
$ws.Range("P8").Value = 0.2120453146491552
$ws.Range("Q8").Value = 148.162522899916
$ws.Range("R8").Value = 1333.462706099244
$ws.Range("S8").Value = 0.01643588083546056
$ws.Range("T8").Value = 0.01643588083546057

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.960212333333333
$ws.Range("H9").Value = 26.880637
$ws.Range("I9").Value = 0.077511171905189
$ws.Range("J9").Value = 0.07751117190518901
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 40.62063066666667
$ws.Range("N9").Value = 121.861892
$ws.Range("O9").Value = 0.5209011059384622
$ws.Range("P9").Value = 0.5209011059384622
$ws.Range("Q9").Value = 363.9694758872449
$ws.Range("R9").Value = 3275.725282985204
$ws.Range("S9").Value = 0.04037565516799921
$ws.Range("T9").Value = 0.04037565516799922

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.960212333333333
$ws.Range("H10").Value = 26.880637
$ws.Range("I10").Value = 0.077511171905189
$ws.Range("J10").Value = 0.07751117190518901
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 20.825229
$ws.Range("N10").Value = 62.475687
$ws.Range("O10").Value = 0.2670535794123827
$ws.Range("P10").Value = 0.2670535794123827
$ws.Range("Q10").Value = 186.598473730291
$ws.Range("R10").Value = 1679.386263572619
$ws.Range("S10").Value = 0.02069963590172923
$ws.Range("T10").Value = 0.02069963590172924
